# Auto-generated Excel COM-interop script applying the market-price refresh diff
# to Marilith_Profits.xlsx (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 6750  # H64: 8500 -> 6750
$ws.Cells.Item(64, 9).Value = 5000  # I64: 0 -> 5000
$ws.Cells.Item(64, 11).Value = 5000  # K64: 0 -> 5000
$ws.Cells.Item(64, 13).Value = -4752  # M64: NEW -> -4752

$ws.Cells.Item(67, 8).Value = 6750  # H67: 8500 -> 6750
$ws.Cells.Item(67, 9).Value = 5000  # I67: 0 -> 5000
$ws.Cells.Item(67, 11).Value = 5000  # K67: 0 -> 5000
$ws.Cells.Item(67, 13).Value = -4142  # M67: NEW -> -4142

$ws.Cells.Item(70, 9).Value = 1351777.8  # I70: 1448190.6 -> 1351777.8
$ws.Cells.Item(70, 10).Value = 3878.6667  # J70: 3761.1875 -> 3878.6667
$ws.Cells.Item(70, 11).Value = 4055333.4  # K70: 4344571.800000001 -> 4055333.4
$ws.Cells.Item(70, 12).Value = 11636.0001  # L70: 11283.5625 -> 11636.0001
$ws.Cells.Item(70, 13).Value = -4055063.4  # M70: -4344301.800000001 -> -4055063.4
$ws.Cells.Item(70, 14).Value = -12176.0001  # N70: -11823.5625 -> -12176.0001

$ws.Cells.Item(73, 9).Value = 1351777.8  # I73: 1448190.6 -> 1351777.8
$ws.Cells.Item(73, 10).Value = 3878.6667  # J73: 3761.1875 -> 3878.6667
$ws.Cells.Item(73, 11).Value = 4055333.4  # K73: 4344571.800000001 -> 4055333.4
$ws.Cells.Item(73, 12).Value = 11636.0001  # L73: 11283.5625 -> 11636.0001
$ws.Cells.Item(73, 13).Value = -4054397.4  # M73: -4343635.800000001 -> -4054397.4
$ws.Cells.Item(73, 14).Value = -13508.0001  # N73: -13155.5625 -> -13508.0001

$ws.Cells.Item(76, 8).Value = 7500  # H76: 8250 -> 7500
$ws.Cells.Item(76, 9).Value = 6000  # I76: 0 -> 6000
$ws.Cells.Item(76, 11).Value = 6000  # K76: 0 -> 6000
$ws.Cells.Item(76, 13).Value = -5685  # M76: NEW -> -5685

$ws.Cells.Item(79, 8).Value = 7500  # H79: 8250 -> 7500
$ws.Cells.Item(79, 9).Value = 6000  # I79: 0 -> 6000
$ws.Cells.Item(79, 11).Value = 6000  # K79: 0 -> 6000
$ws.Cells.Item(79, 13).Value = -4908  # M79: NEW -> -4908

$ws.Cells.Item(93, 8).Value = 44119.6  # H93: 46239.6 -> 44119.6
$ws.Cells.Item(93, 10).Value = 44119.6  # J93: 46239.6 -> 44119.6
$ws.Cells.Item(93, 12).Value = 44119.6  # L93: 46239.6 -> 44119.6
$ws.Cells.Item(93, 14).Value = -49111.6  # N93: -51231.6 -> -49111.6

$ws.Cells.Item(125, 8).Value = 8499.75  # H125: 7166.5 -> 8499.75
$ws.Cells.Item(125, 9).Value = 5000  # I125: 4750 -> 5000
$ws.Cells.Item(125, 11).Value = 45000  # K125: 42750 -> 45000
$ws.Cells.Item(125, 13).Value = -42540  # M125: -40290 -> -42540

$ws.Cells.Item(132, 8).Value = 26318844  # H132: 27030148 -> 26318844
$ws.Cells.Item(132, 9).Value = 31252972  # I132: 32261114 -> 31252972
$ws.Cells.Item(132, 11).Value = 93758916  # K132: 96783342 -> 93758916
$ws.Cells.Item(132, 13).Value = -93756386  # M132: -96780812 -> -93756386

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 8156.027  # H32: 8104.9116 -> 8156.027
$ws.Cells.Item(32, 9).Value = 8156.027  # I32: 8199.182000000001 -> 8156.027
$ws.Cells.Item(32, 10).Value = 0  # J32: 4994 -> 0
$ws.Cells.Item(32, 11).Value = 8156.027  # K32: 8199.182000000001 -> 8156.027
$ws.Cells.Item(32, 12).Value = 0  # L32: 4994 -> 0
$ws.Cells.Item(32, 13).Value = -7869.027  # M32: -7912.182000000001 -> -7869.027
$ws.Cells.Item(32, 14).ClearContents()  # N32: -5568 -> (removed)

$ws.Cells.Item(61, 8).Value = 2300  # H61: 2000 -> 2300
$ws.Cells.Item(61, 10).Value = 2900  # J61: 0 -> 2900
$ws.Cells.Item(61, 12).Value = 2900  # L61: 0 -> 2900
$ws.Cells.Item(61, 14).Value = -3324  # N61: NEW -> -3324

$ws.Cells.Item(63, 8).Value = 2605.111  # H63: 2394.5 -> 2605.111
$ws.Cells.Item(63, 9).Value = 2089.4  # I63: 1824.3334 -> 2089.4
$ws.Cells.Item(63, 11).Value = 2089.4  # K63: 1824.3334 -> 2089.4
$ws.Cells.Item(63, 13).Value = -1403.4  # M63: -1138.3334 -> -1403.4

$ws.Cells.Item(66, 8).Value = 2605.111  # H66: 2394.5 -> 2605.111
$ws.Cells.Item(66, 9).Value = 2089.4  # I66: 1824.3334 -> 2089.4
$ws.Cells.Item(66, 11).Value = 10447  # K66: 9121.666999999999 -> 10447
$ws.Cells.Item(66, 13).Value = -7015  # M66: -5689.666999999999 -> -7015

$ws.Cells.Item(96, 8).Value = 13921.75  # H96: 12137.4 -> 13921.75
$ws.Cells.Item(96, 10).Value = 13921.75  # J96: 12137.4 -> 13921.75
$ws.Cells.Item(96, 12).Value = 13921.75  # L96: 12137.4 -> 13921.75
$ws.Cells.Item(96, 14).Value = -19413.75  # N96: -17629.4 -> -19413.75

$ws.Cells.Item(122, 8).Value = 3599.0715  # H122: 3798.9167 -> 3599.0715
$ws.Cells.Item(122, 9).Value = 3645.1538  # I122: 3871.5454 -> 3645.1538
$ws.Cells.Item(122, 11).Value = 10935.4614  # K122: 11614.6362 -> 10935.4614
$ws.Cells.Item(122, 13).Value = -8485.4614  # M122: -9164.636200000001 -> -8485.4614

$ws.Cells.Item(132, 8).Value = 4099.2856  # H132: 4022.0454 -> 4099.2856
$ws.Cells.Item(132, 9).Value = 4099.2856  # I132: 4022.0454 -> 4099.2856
$ws.Cells.Item(132, 11).Value = 12297.8568  # K132: 12066.1362 -> 12297.8568
$ws.Cells.Item(132, 13).Value = -9767.856800000001  # M132: -9536.136200000001 -> -9767.856800000001

$ws.Cells.Item(136, 8).Value = 2300  # H136: 2000 -> 2300
$ws.Cells.Item(136, 10).Value = 2900  # J136: 0 -> 2900
$ws.Cells.Item(136, 12).Value = 8700  # L136: 0 -> 8700
$ws.Cells.Item(136, 14).Value = -13800  # N136: NEW -> -13800

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(12, 8).Value = 203.16667  # H12: 1126.25 -> 203.16667
$ws.Cells.Item(12, 9).Value = 174.66667  # I12: 102.5 -> 174.66667
$ws.Cells.Item(12, 10).Value = 231.66667  # J12: 2150 -> 231.66667
$ws.Cells.Item(12, 11).Value = 174.66667  # K12: 102.5 -> 174.66667
$ws.Cells.Item(12, 12).Value = 231.66667  # L12: 2150 -> 231.66667
$ws.Cells.Item(12, 13).Value = -6.666670000000011  # M12: 65.5 -> -6.666670000000011
$ws.Cells.Item(12, 14).Value = -567.6666700000001  # N12: -2486 -> -567.6666700000001

$ws.Cells.Item(94, 8).Value = 2614.8462  # H94: 2676.5833 -> 2614.8462
$ws.Cells.Item(94, 9).Value = 2624.4167  # I94: 2692.6365 -> 2624.4167
$ws.Cells.Item(94, 11).Value = 2624.4167  # K94: 2692.6365 -> 2624.4167
$ws.Cells.Item(94, 13).Value = -2173.4167  # M94: -2241.6365 -> -2173.4167

$ws.Cells.Item(105, 8).Value = 2577.8572  # H105: 2358.5715 -> 2577.8572
$ws.Cells.Item(105, 9).Value = 1010  # I105: 920 -> 1010
$ws.Cells.Item(105, 10).Value = 2839.1667  # J105: 3437.5 -> 2839.1667
$ws.Cells.Item(105, 11).Value = 1010  # K105: 920 -> 1010
$ws.Cells.Item(105, 12).Value = 2839.1667  # L105: 3437.5 -> 2839.1667
$ws.Cells.Item(105, 13).Value = 737  # M105: 827 -> 737
$ws.Cells.Item(105, 14).Value = -6333.1667  # N105: -6931.5 -> -6333.1667

$ws.Cells.Item(134, 8).Value = 5524.294  # H134: 5521.353 -> 5524.294
$ws.Cells.Item(134, 9).Value = 5407.6  # I134: 5404.2666 -> 5407.6
$ws.Cells.Item(134, 11).Value = 16222.8  # K134: 16212.7998 -> 16222.8
$ws.Cells.Item(134, 13).Value = -13687.8  # M134: -13677.7998 -> -13687.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 587.5  # H22: 650 -> 587.5
$ws.Cells.Item(22, 9).Value = 575  # I22: 650 -> 575
$ws.Cells.Item(22, 10).Value = 600  # J22: 0 -> 600
$ws.Cells.Item(22, 11).Value = 575  # K22: 650 -> 575
$ws.Cells.Item(22, 12).Value = 600  # L22: 0 -> 600
$ws.Cells.Item(22, 13).Value = -225  # M22: -300 -> -225
$ws.Cells.Item(22, 14).Value = -1300  # N22: NEW -> -1300

$ws.Cells.Item(31, 8).Value = 4059.8125  # H31: 4026.9412 -> 4059.8125
$ws.Cells.Item(31, 9).Value = 4499.5  # I31: 4299.6 -> 4499.5
$ws.Cells.Item(31, 10).Value = 3913.25  # J31: 3913.3333 -> 3913.25
$ws.Cells.Item(31, 11).Value = 4499.5  # K31: 4299.6 -> 4499.5
$ws.Cells.Item(31, 12).Value = 3913.25  # L31: 3913.3333 -> 3913.25
$ws.Cells.Item(31, 13).Value = -4204.5  # M31: -4004.6 -> -4204.5
$ws.Cells.Item(31, 14).Value = -4503.25  # N31: -4503.3333 -> -4503.25

$ws.Cells.Item(34, 8).Value = 4059.8125  # H34: 4026.9412 -> 4059.8125
$ws.Cells.Item(34, 9).Value = 4499.5  # I34: 4299.6 -> 4499.5
$ws.Cells.Item(34, 10).Value = 3913.25  # J34: 3913.3333 -> 3913.25
$ws.Cells.Item(34, 11).Value = 4499.5  # K34: 4299.6 -> 4499.5
$ws.Cells.Item(34, 12).Value = 3913.25  # L34: 3913.3333 -> 3913.25
$ws.Cells.Item(34, 13).Value = -4297.5  # M34: -4097.6 -> -4297.5
$ws.Cells.Item(34, 14).Value = -4317.25  # N34: -4317.3333 -> -4317.25

$ws.Cells.Item(57, 8).Value = 16000  # H57: 0 -> 16000
$ws.Cells.Item(57, 10).Value = 16000  # J57: 0 -> 16000
$ws.Cells.Item(57, 12).Value = 16000  # L57: 0 -> 16000
$ws.Cells.Item(57, 14).Value = -17120  # N57: NEW -> -17120

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(28, 8).Value = 3865  # H28: 0 -> 3865
$ws.Cells.Item(28, 9).Value = 2730  # I28: 0 -> 2730
$ws.Cells.Item(28, 10).Value = 5000  # J28: 0 -> 5000
$ws.Cells.Item(28, 11).Value = 8190  # K28: 0 -> 8190
$ws.Cells.Item(28, 12).Value = 15000  # L28: 0 -> 15000
$ws.Cells.Item(28, 13).Value = -7958  # M28: NEW -> -7958
$ws.Cells.Item(28, 14).Value = -15464  # N28: NEW -> -15464

$ws.Cells.Item(29, 8).Value = 688.1429000000001  # H29: 714 -> 688.1429000000001
$ws.Cells.Item(29, 9).Value = 189.5  # I29: 280 -> 189.5
$ws.Cells.Item(29, 11).Value = 568.5  # K29: 840 -> 568.5
$ws.Cells.Item(29, 13).Value = -291.5  # M29: -563 -> -291.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 7869  # H70: 7651.75 -> 7869
$ws.Cells.Item(70, 10).Value = 0  # J70: 7000 -> 0
$ws.Cells.Item(70, 12).Value = 0  # L70: 7000 -> 0
$ws.Cells.Item(70, 14).ClearContents()  # N70: -7540 -> (removed)

$ws.Cells.Item(73, 8).Value = 7869  # H73: 7651.75 -> 7869
$ws.Cells.Item(73, 10).Value = 0  # J73: 7000 -> 0
$ws.Cells.Item(73, 12).Value = 0  # L73: 7000 -> 0
$ws.Cells.Item(73, 14).ClearContents()  # N73: -8872 -> (removed)

$ws.Cells.Item(80, 8).Value = 6916.3335  # H80: 6214 -> 6916.3335
$ws.Cells.Item(80, 9).Value = 6249.5  # I80: 4833 -> 6249.5
$ws.Cells.Item(80, 11).Value = 6249.5  # K80: 4833 -> 6249.5
$ws.Cells.Item(80, 13).Value = -5251.5  # M80: -3835 -> -5251.5

$ws.Cells.Item(83, 8).Value = 6916.3335  # H83: 6214 -> 6916.3335
$ws.Cells.Item(83, 9).Value = 6249.5  # I83: 4833 -> 6249.5
$ws.Cells.Item(83, 11).Value = 31247.5  # K83: 24165 -> 31247.5
$ws.Cells.Item(83, 13).Value = -26255.5  # M83: -19173 -> -26255.5

$ws.Cells.Item(98, 8).Value = 20000  # H98: 10643 -> 20000
$ws.Cells.Item(98, 10).Value = 20000  # J98: 10643 -> 20000
$ws.Cells.Item(98, 12).Value = 20000  # L98: 10643 -> 20000
$ws.Cells.Item(98, 14).Value = -25990  # N98: -16633 -> -25990

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1292.619  # H22: 1292.5238 -> 1292.619
$ws.Cells.Item(22, 9).Value = 1140  # I22: 1130.4667 -> 1140
$ws.Cells.Item(22, 10).Value = 1597.8572  # J22: 1697.6666 -> 1597.8572
$ws.Cells.Item(22, 11).Value = 1140  # K22: 1130.4667 -> 1140
$ws.Cells.Item(22, 12).Value = 1597.8572  # L22: 1697.6666 -> 1597.8572
$ws.Cells.Item(22, 13).Value = -845  # M22: -835.4666999999999 -> -845
$ws.Cells.Item(22, 14).Value = -2187.8572  # N22: -2287.6666 -> -2187.8572

$ws.Cells.Item(27, 8).Value = 1292.619  # H27: 1292.5238 -> 1292.619
$ws.Cells.Item(27, 9).Value = 1140  # I27: 1130.4667 -> 1140
$ws.Cells.Item(27, 10).Value = 1597.8572  # J27: 1697.6666 -> 1597.8572
$ws.Cells.Item(27, 11).Value = 1140  # K27: 1130.4667 -> 1140
$ws.Cells.Item(27, 12).Value = 1597.8572  # L27: 1697.6666 -> 1597.8572
$ws.Cells.Item(27, 13).Value = -1033  # M27: -1023.4667 -> -1033
$ws.Cells.Item(27, 14).Value = -1811.8572  # N27: -1911.6666 -> -1811.8572

$ws.Cells.Item(104, 8).Value = 17456.666  # H104: 21342.5 -> 17456.666
$ws.Cells.Item(104, 10).Value = 17456.666  # J104: 21342.5 -> 17456.666
$ws.Cells.Item(104, 12).Value = 17456.666  # L104: 21342.5 -> 17456.666
$ws.Cells.Item(104, 14).Value = -24444.666  # N104: -28330.5 -> -24444.666

$ws.Cells.Item(132, 8).Value = 26497.75  # H132: 28775.889 -> 26497.75
$ws.Cells.Item(132, 9).Value = 25270.727  # I132: 27373.5 -> 25270.727
$ws.Cells.Item(132, 11).Value = 75812.181  # K132: 82120.5 -> 75812.181
$ws.Cells.Item(132, 13).Value = -73282.181  # M132: -79590.5 -> -73282.181

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(45, 8).Value = 0  # H45: 14425 -> 0
$ws.Cells.Item(45, 9).Value = 0  # I45: 14425 -> 0
$ws.Cells.Item(45, 11).Value = 0  # K45: 14425 -> 0
$ws.Cells.Item(45, 13).ClearContents()  # M45: -13934 -> (removed)

$ws.Cells.Item(122, 8).Value = 1202.3846  # H122: 1147.5714 -> 1202.3846
$ws.Cells.Item(122, 9).Value = 1202.3846  # I122: 1147.5714 -> 1202.3846
$ws.Cells.Item(122, 11).Value = 3607.1538  # K122: 3442.7142 -> 3607.1538
$ws.Cells.Item(122, 13).Value = -1157.1538  # M122: -992.7142000000003 -> -1157.1538

$ws.Cells.Item(140, 8).Value = 68666.336  # H140: 60000 -> 68666.336
$ws.Cells.Item(140, 10).Value = 68666.336  # J140: 60000 -> 68666.336
$ws.Cells.Item(140, 12).Value = 68666.336  # L140: 60000 -> 68666.336
$ws.Cells.Item(140, 14).Value = -79026.336  # N140: -70360 -> -79026.336
